{"js": "// Update the date paragraph and the 25 \"problem\" table cells to the new\n// values described by the commit (two-digit \u00f7 one-digit division drills).\n//\n// The table has 5 columns; every 4th row (0, 4, 8, 12, 16) holds the five\n// visible problems for that block (the 3 rows following each are blank\n// workspace rows). Two of the ORIGINAL cells share identical text\n// (\"92\u00f75=18, 2\" appears twice in row 4), so a naive global find/replace\n// would not be able to give them different new values \u2014 addressing cells\n// by (row, column) sidesteps that ambiguity entirely.\n\nconst body = context.document.body;\n\n// --- 1. Update the date heading paragraph -------------------------------\nconst dateParas = body.paragraphs.search(\"2024-07-20 Saturday\", { matchCase: true });\ndateParas.load(\"items\");\nawait context.sync();\nif (dateParas.items.length > 0) {\n  dateParas.items[0].insertText(\"2024-07-21 Sunday\", \"Replace\");\n} else {\n  // Fallback: first paragraph of the document holds the date.\n  const firstPara = body.paragraphs.getFirst();\n  firstPara.insertText(\"2024-07-21 Sunday\", \"Replace\");\n}\nawait context.sync();\n\n// --- 2. Update the table cells ------------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, colIndex, newText]\nconst updates = [\n  [0, 0, \"39\u00f75=7, 4\"],\n  [0, 1, \"19\u00f73=6, 1\"],\n  [0, 2, \"50\u00f78=6, 2\"],\n  [0, 3, \"52\u00f78=6, 4\"],\n  [0, 4, \"70\u00f76=11, 4\"],\n  [4, 0, \"18\u00f77=2, 4\"],\n  [4, 1, \"54\u00f76=9, 0\"],\n  [4, 2, \"29\u00f77=4, 1\"],\n  [4, 3, \"41\u00f72=20, 1\"],\n  [4, 4, \"74\u00f73=24, 2\"],\n  [8, 0, \"30\u00f79=3, 3\"],\n  [8, 1, \"53\u00f75=10, 3\"],\n  [8, 2, \"22\u00f75=4, 2\"],\n  [8, 3, \"61\u00f76=10, 1\"],\n  [8, 4, \"41\u00f73=13, 2\"],\n  [12, 0, \"71\u00f74=17, 3\"],\n  [12, 1, \"17\u00f73=5, 2\"],\n  [12, 2, \"48\u00f78=6, 0\"],\n  [12, 3, \"59\u00f72=29, 1\"],\n  [12, 4, \"27\u00f72=13, 1\"],\n  [16, 0, \"21\u00f78=2, 5\"],\n  [16, 1, \"93\u00f79=10, 3\"],\n  [16, 2, \"52\u00f72=26, 0\"],\n  [16, 3, \"44\u00f75=8, 4\"],\n  [16, 4, \"30\u00f72=15, 0\"],\n];\n\nfor (const [row, col, text] of updates) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the date paragraph and the 25 \"problem\" table cells to the new\n# values described by the commit (two-digit \u00f7 one-digit division drills).\n#\n# The table has 5 columns; every 4th row (1, 5, 9, 13, 17 in 1-based COM\n# indexing) holds the five visible problems for that block (the 3 rows\n# following each are blank workspace rows). Two of the ORIGINAL cells share\n# identical text (\"92\u00f75=18, 2\" appears twice in row 5), so a naive global\n# Find/Replace could not give them two different new values -- addressing\n# cells by (row, column) via Table.Cell() sidesteps that ambiguity entirely.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date heading paragraph --------------------------------\n$find = $d.Content.Find\n$find.Text = \"2024-07-20 Saturday\"\n$find.Replacement.Text = \"2024-07-21 Sunday\"\n$find.Execute(\n    [ref]$find.Text,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]1,\n    [ref]$false,\n    [ref]$find.Replacement.Text,\n    [ref]2\n)\n\n# --- 2. Update the table cells --------------------------------------------\n$t = $d.Tables.Item(1)\n\n# (row, col, newText) -- 1-based row/col, matching Word COM Table.Cell(row, col)\n$updates = @(\n    @(1, 1, \"39\u00f75=7, 4\"),\n    @(1, 2, \"19\u00f73=6, 1\"),\n    @(1, 3, \"50\u00f78=6, 2\"),\n    @(1, 4, \"52\u00f78=6, 4\"),\n    @(1, 5, \"70\u00f76=11, 4\"),\n    @(5, 1, \"18\u00f77=2, 4\"),\n    @(5, 2, \"54\u00f76=9, 0\"),\n    @(5, 3, \"29\u00f77=4, 1\"),\n    @(5, 4, \"41\u00f72=20, 1\"),\n    @(5, 5, \"74\u00f73=24, 2\"),\n    @(9, 1, \"30\u00f79=3, 3\"),\n    @(9, 2, \"53\u00f75=10, 3\"),\n    @(9, 3, \"22\u00f75=4, 2\"),\n    @(9, 4, \"61\u00f76=10, 1\"),\n    @(9, 5, \"41\u00f73=13, 2\"),\n    @(13, 1, \"71\u00f74=17, 3\"),\n    @(13, 2, \"17\u00f73=5, 2\"),\n    @(13, 3, \"48\u00f78=6, 0\"),\n    @(13, 4, \"59\u00f72=29, 1\"),\n    @(13, 5, \"27\u00f72=13, 1\"),\n    @(17, 1, \"21\u00f78=2, 5\"),\n    @(17, 2, \"93\u00f79=10, 3\"),\n    @(17, 3, \"52\u00f72=26, 0\"),\n    @(17, 4, \"44\u00f75=8, 4\"),\n    @(17, 5, \"30\u00f72=15, 0\")\n)\n\nforeach ($u in $updates) {\n    $cell = $t.Cell($u[0], $u[1])\n    $cell.Range.Text = $u[2]\n}\n"}
